# Script 1 - atualização automática de dados (Execução: 23)
# Updates the "g1.1" sheet: refreshes values/ranks for "Variação (%) 2022" ->
# "Variação (%) 2023" and "Variação (%) 2022/2010" -> "Variação (%) 2023/2010"
# blocks with the newly computed figures (states get re-ranked too).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-10: "Variação (%) 2023" block
$ws.Cells.Item(2,1).Value = "AC"
$ws.Cells.Item(2,2).Value = 14.73257689442189
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(2,4).Value = "Variação (%) 2023"

$ws.Cells.Item(3,1).Value = "MS"
$ws.Cells.Item(3,2).Value = 13.44269577606423
$ws.Cells.Item(3,3).Value = 2
$ws.Cells.Item(3,4).Value = "Variação (%) 2023"

$ws.Cells.Item(4,1).Value = "MT"
$ws.Cells.Item(4,2).Value = 12.88001598426398
$ws.Cells.Item(4,3).Value = 3
$ws.Cells.Item(4,4).Value = "Variação (%) 2023"

$ws.Cells.Item(5,1).Value = "TO"
$ws.Cells.Item(5,2).Value = 7.890383025089162
$ws.Cells.Item(5,3).Value = 4
$ws.Cells.Item(5,4).Value = "Variação (%) 2023"

$ws.Cells.Item(6,1).Value = "RJ"
$ws.Cells.Item(6,2).Value = 5.652659822157795
$ws.Cells.Item(6,3).Value = 5
$ws.Cells.Item(6,4).Value = "Variação (%) 2023"

$ws.Cells.Item(7,1).Value = "GO"
$ws.Cells.Item(7,2).Value = 4.816953216278661
$ws.Cells.Item(7,3).Value = 6
$ws.Cells.Item(7,4).Value = "Variação (%) 2023"

$ws.Cells.Item(8,1).Value = "SE"
$ws.Cells.Item(8,2).Value = 3.118144130554446
$ws.Cells.Item(8,3).Value = 15
$ws.Cells.Item(8,4).Value = "Variação (%) 2023"

$ws.Cells.Item(9,1).Value = "BR"
$ws.Cells.Item(9,2).Value = 3.241657824791806
$ws.Cells.Item(9,4).Value = "Variação (%) 2023"

$ws.Cells.Item(10,1).Value = "NE"
$ws.Cells.Item(10,2).Value = 2.867008788862638
$ws.Cells.Item(10,4).Value = "Variação (%) 2023"

# Row 11-19: "Variação (%) 2023/2010" block
$ws.Cells.Item(11,1).Value = "MT"
$ws.Cells.Item(11,2).Value = 77.14346626765018
$ws.Cells.Item(11,3).Value = 1
$ws.Cells.Item(11,4).Value = "Variação (%) 2023/2010"

$ws.Cells.Item(12,1).Value = "RR"
$ws.Cells.Item(12,2).Value = 63.84407511155798
$ws.Cells.Item(12,3).Value = 2
$ws.Cells.Item(12,4).Value = "Variação (%) 2023/2010"

$ws.Cells.Item(13,1).Value = "TO"
$ws.Cells.Item(13,2).Value = 59.35173933449352
$ws.Cells.Item(13,3).Value = 3
$ws.Cells.Item(13,4).Value = "Variação (%) 2023/2010"

$ws.Cells.Item(14,1).Value = "MS"
$ws.Cells.Item(14,2).Value = 49.57179111911111
$ws.Cells.Item(14,3).Value = 4
$ws.Cells.Item(14,4).Value = "Variação (%) 2023/2010"

$ws.Cells.Item(15,1).Value = "AC"
$ws.Cells.Item(15,2).Value = 42.69524774665621
$ws.Cells.Item(15,3).Value = 5
$ws.Cells.Item(15,4).Value = "Variação (%) 2023/2010"

$ws.Cells.Item(16,1).Value = "PI"
$ws.Cells.Item(16,2).Value = 36.70681689547283
$ws.Cells.Item(16,3).Value = 6
$ws.Cells.Item(16,4).Value = "Variação (%) 2023/2010"

$ws.Cells.Item(17,1).Value = "SE"
$ws.Cells.Item(17,2).Value = 7.324239245718005
$ws.Cells.Item(17,3).Value = 25
$ws.Cells.Item(17,4).Value = "Variação (%) 2023/2010"

$ws.Cells.Item(18,1).Value = "BR"
$ws.Cells.Item(18,2).Value = 15.14430956101356
$ws.Cells.Item(18,4).Value = "Variação (%) 2023/2010"

$ws.Cells.Item(19,1).Value = "NE"
$ws.Cells.Item(19,2).Value = 16.86384673819174
$ws.Cells.Item(19,4).Value = "Variação (%) 2023/2010"
